# Update the "as of" date in the confidential disclaimer text (shared string
# used by cell A18) from 2021-04-29 to 2021-04-30.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected (no formatting/structure changes are being made,
# only cell values), so unprotect it first to allow writing the cells, then
# restore protection afterwards.
$ws.Unprotect()

$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-30 for illustrative purposes only and are subject to change."
# Assigning a value with an embedded line break makes the host resize row 18
# to fit two lines; re-run AutoFit so the row reverts to the sheet's default
# (un-customized) height, matching the rest of the file.
$ws.Rows.Item(18).AutoFit()

# Updated Weight (D) and Percent Change (E) values for rows 2-15.
$ws.Range("D2").Value  = 0.05794108462027825
$ws.Range("E2").Value  = -0.007315957933241912

$ws.Range("D3").Value  = 0.02378909392211401
$ws.Range("E3").Value  = -0.01364256480218295

$ws.Range("D4").Value  = 0.03205238806617656
$ws.Range("E4").Value  = -0.01469220755067879

$ws.Range("D5").Value  = 0.03185103816077459
$ws.Range("E5").Value  = -0.02331094429079428

$ws.Range("D6").Value  = 0.0373159660208788
$ws.Range("E6").Value  = -0.02428115015974452

$ws.Range("D7").Value  = 0.01919679270795422
$ws.Range("E7").Value  = -0.0139983093829249

$ws.Range("D8").Value  = 0.004703242362695426
$ws.Range("E8").Value  = -0.0181664554288129

$ws.Range("D9").Value  = 0.006962467781531379
$ws.Range("E9").Value  = -0.007229832572298367

$ws.Range("D10").Value = 0.07006181905730491
$ws.Range("E10").Value = -0.002836074872376582

$ws.Range("D11").Value = 0.07018103939602975
$ws.Range("E11").Value = -0.003397508493771184

$ws.Range("D12").Value = 0.1465827311326363
$ws.Range("E12").Value = 0.002313475997686387

$ws.Range("D13").Value = 0.3847101240255807
$ws.Range("E13").Value = 0.0008780402142416488

$ws.Range("D14").Value = 0.1146522127460451
$ws.Range("E14").Value = -0.004488657042339028

$ws.Range("E15").Value = -0.003547281958419357

# Restore sheet protection to match the original protected state.
$ws.Protect()
